$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Raw parking data: park ID, vehicle id, slot ID, in time (as text)
# -----------------------------------------------------------------
$parkId  = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$vehId   = @(4, 1, 3, 8, 12, 17, 6, 19, 7, 11)
$slotId  = @(2, 65, 45, 87, 112, 118, 194, 154, 100, 34)
$inTime  = @(
    "2018-06-03 06:00:00",
    "2018-06-03 06:30:00",
    "2018-06-03 07:00:00",
    "2018-06-03 07:10:00",
    "2018-06-03 07:40:00",
    "2018-06-03 07:50:00",
    "2018-06-03 08:30:00",
    "2018-06-13 09:30:00",
    "2018-06-13 10:00:00",
    "2018-06-13 11:00:00"
)

$n = $parkId.Length

# Columns A, B, C: plain numbers
for ($i = 0; $i -lt $n; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $parkId[$i]
    $ws.Cells.Item($r, 2).Value = $vehId[$i]
    $ws.Cells.Item($r, 3).Value = $slotId[$i]
}

# Column D: "in time" stored as TEXT (number format "@" == numFmtId 49)
$dRange = $ws.Range("D1:D$n")
$dRange.NumberFormat = "@"
for ($i = 0; $i -lt $n; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 4).Value = $inTime[$i]
}

# Column E: CONCATENATE formula building the SQL insert statement.
# E1 gets its own (non-shared) formula; E2:E10 are entered as one range
# assignment so Excel stores them as a shared formula group.
$ws.Range("E1").Formula = '=CONCATENATE("insert into [dbo].[parking] ([park ID],[vehicle id],[slot ID],[in time]) values(",A1,",",B1,",",C1,",' + "'" + '",D1,"' + "'" + ')")'
$ws.Range("E2:E$n").Formula = '=CONCATENATE("insert into [dbo].[parking] ([park ID],[vehicle id],[slot ID],[in time]) values(",A2,",",B2,",",C2,",' + "'" + '",D2,"' + "'" + ')")'

# Column widths (D, E, F) -- chosen so the character width converts (via the
# host's pixel-rounding) to the closest achievable value to the widths
# recorded in the source file (19, 99.5703125, 9.28515625)
$ws.Columns.Item(4).ColumnWidth = 18.14
$ws.Columns.Item(5).ColumnWidth = 98.65
$ws.Columns.Item(6).ColumnWidth = 8.5

# Row 7 is slightly taller in the source file
$ws.Rows.Item(7).RowHeight = 16.5

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Selection / active cell as left by the author
$ws.Range("E1:E$n").Select()
